# Regenerate save_data: update the "K" column (column G) values for each
# row with the newly computed strike counts (K instead of Strike#).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 5
    3  = 4
    4  = 3
    5  = 11
    6  = 5
    7  = 9
    8  = 9
    9  = 2
    10 = 4
    11 = 5
    12 = 2
    13 = 2
    14 = 6
    15 = 8
    16 = 4
    17 = 7
    18 = 1
    19 = 3
    20 = 7
    21 = 6
    22 = 6
    23 = 3
    24 = 7
    25 = 5
    26 = 4
    27 = 9
    28 = 9
    29 = 6
    30 = 5
    31 = 7
    32 = 4
    33 = 4
    34 = 7
    35 = 4
    36 = 4
    37 = 6
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newValues[$row]
}
